# Auto-applied update of market-price derived columns (H:N) across all Sheets
# per scheduled market-data refresh. Values taken verbatim from the refreshed
# dataset; some rows gain/lose Profit (M/N) cells depending on whether a craft
# is currently profitable.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2321.5
$ws.Range("I2").Value = 630
$ws.Range("J2").Value = 2659.8
$ws.Range("K2").Value = 630
$ws.Range("L2").Value = 2659.8
$ws.Range("M2").Value = -517
$ws.Range("N2").Value = -2885.8
$ws.Range("H4").Value = 117
$ws.Range("I4").Value = 130.16667
$ws.Range("J4").Value = 77.5
$ws.Range("K4").Value = 130.16667
$ws.Range("L4").Value = 77.5
$ws.Range("M4").Value = -16.16667000000001
$ws.Range("N4").Value = -305.5
$ws.Range("H32").Value = 941.4286
$ws.Range("I32").Value = 995
$ws.Range("J32").Value = 932.5
$ws.Range("K32").Value = 995
$ws.Range("L32").Value = 932.5
$ws.Range("M32").Value = -669
$ws.Range("N32").Value = -1584.5
$ws.Range("H64").Value = 4530.2
$ws.Range("I64").Value = 4530.2
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 4530.2
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -4282.2
$ws.Range("H67").Value = 4530.2
$ws.Range("I67").Value = 4530.2
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 4530.2
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -3672.2
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("H100").Value = 3100
$ws.Range("I100").Value = 3050
$ws.Range("J100").Value = 3300
$ws.Range("K100").Value = 3050
$ws.Range("L100").Value = 3300
$ws.Range("M100").Value = -2509
$ws.Range("N100").Value = -4382
$ws.Range("H132").Value = 3252
$ws.Range("I132").Value = 3252
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9756
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7226
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("H137").Value = 1888.3636
$ws.Range("I137").Value = 977
$ws.Range("J137").Value = 2230.125
$ws.Range("K137").Value = 2931
$ws.Range("L137").Value = 6690.375
$ws.Range("M137").Value = -381
$ws.Range("N137").Value = -11790.375
$ws.Range("M135").ClearContents()
$ws.Range("M69").ClearContents()
$ws.Range("M72").ClearContents()
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1546.6562
$ws.Range("I32").Value = 1570.742
$ws.Range("J32").Value = 800
$ws.Range("K32").Value = 1570.742
$ws.Range("L32").Value = 800
$ws.Range("M32").Value = -1283.742
$ws.Range("H74").Value = 4043.9
$ws.Range("I74").Value = 4043.9
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 4043.9
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -3169.9
$ws.Range("H77").Value = 4043.9
$ws.Range("I77").Value = 4043.9
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 20219.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -15851.5
$ws.Range("H102").Value = 1872.5
$ws.Range("I102").Value = 1654.2858
$ws.Range("J102").Value = 3400
$ws.Range("K102").Value = 1654.2858
$ws.Range("L102").Value = 3400
$ws.Range("M102").Value = -32.28580000000011
$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6550
$ws.Range("H132").Value = 1436.24
$ws.Range("I132").Value = 1364.619
$ws.Range("J132").Value = 1812.25
$ws.Range("K132").Value = 4093.857
$ws.Range("L132").Value = 5436.75
$ws.Range("M132").Value = -1563.857
$ws.Range("N132").Value = -10496.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2586.2222
$ws.Range("I94").Value = 2586.2222
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2586.2222
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -2135.2222
$ws.Range("H107").Value = 1065
$ws.Range("I107").Value = 997.5
$ws.Range("J107").Value = 1200
$ws.Range("K107").Value = 997.5
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = 922.5
$ws.Range("H134").Value = 3153.25
$ws.Range("I134").Value = 3349.0908
$ws.Range("J134").Value = 999
$ws.Range("K134").Value = 10047.2724
$ws.Range("L134").Value = 2997
$ws.Range("M134").Value = -7512.2724

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2196.2563
$ws.Range("I31").Value = 1661.9615
$ws.Range("J31").Value = 3264.8462
$ws.Range("K31").Value = 1661.9615
$ws.Range("L31").Value = 3264.8462
$ws.Range("M31").Value = -1366.9615
$ws.Range("N31").Value = -3854.8462
$ws.Range("H34").Value = 2196.2563
$ws.Range("I34").Value = 1661.9615
$ws.Range("J34").Value = 3264.8462
$ws.Range("K34").Value = 1661.9615
$ws.Range("L34").Value = 3264.8462
$ws.Range("M34").Value = -1459.9615
$ws.Range("N34").Value = -3668.8462
$ws.Range("H58").Value = 4234.6113
$ws.Range("I58").Value = 4309.5293
$ws.Range("J58").Value = 2961
$ws.Range("K58").Value = 4309.5293
$ws.Range("L58").Value = 2961
$ws.Range("M58").Value = -4106.5293
$ws.Range("H93").Value = 15476.75
$ws.Range("I93").Value = 15476.75
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 15476.75
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -13604.75
$ws.Range("H132").Value = 3398.3333
$ws.Range("I132").Value = 3395
$ws.Range("J132").Value = 3400
$ws.Range("K132").Value = 10185
$ws.Range("L132").Value = 10200
$ws.Range("M132").Value = -7655
$ws.Range("H134").Value = 3549.8333
$ws.Range("I134").Value = 1650
$ws.Range("J134").Value = 4499.75
$ws.Range("K134").Value = 4950
$ws.Range("L134").Value = 13499.25
$ws.Range("M134").Value = -2415
$ws.Range("N134").Value = -18569.25
$ws.Range("H136").Value = 4234.6113
$ws.Range("I136").Value = 4309.5293
$ws.Range("J136").Value = 2961
$ws.Range("K136").Value = 12928.5879
$ws.Range("L136").Value = 8883
$ws.Range("M136").Value = -10378.5879

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1806.5652
$ws.Range("I4").Value = 1796.8125
$ws.Range("J4").Value = 1828.8572
$ws.Range("K4").Value = 5390.4375
$ws.Range("L4").Value = 5486.571599999999
$ws.Range("M4").Value = -5278.4375
$ws.Range("N4").Value = -5710.571599999999
$ws.Range("H60").Value = 607.8333
$ws.Range("I60").Value = 349.4
$ws.Range("J60").Value = 1900
$ws.Range("K60").Value = 1048.2
$ws.Range("L60").Value = 5700
$ws.Range("M60").Value = -797.1999999999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 2950
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 2950
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 2950
$ws.Range("N7").Value = -3174
$ws.Range("H8").Value = 2950
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 2950
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 2950
$ws.Range("N8").Value = -3228
$ws.Range("H80").Value = 2119.5833
$ws.Range("I80").Value = 1938.7778
$ws.Range("J80").Value = 2662
$ws.Range("K80").Value = 1938.7778
$ws.Range("L80").Value = 2662
$ws.Range("M80").Value = -940.7778000000001
$ws.Range("H83").Value = 2119.5833
$ws.Range("I83").Value = 1938.7778
$ws.Range("J83").Value = 2662
$ws.Range("K83").Value = 9693.889000000001
$ws.Range("L83").Value = 13310
$ws.Range("M83").Value = -4701.889000000001
$ws.Range("H122").Value = 2974.0715
$ws.Range("I122").Value = 1614.3
$ws.Range("J122").Value = 6373.5
$ws.Range("K122").Value = 4842.9
$ws.Range("L122").Value = 19120.5
$ws.Range("M122").Value = -2392.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1220
$ws.Range("I9").Value = 650
$ws.Range("J9").Value = 3500
$ws.Range("K9").Value = 650
$ws.Range("L9").Value = 3500
$ws.Range("M9").Value = -426
$ws.Range("H22").Value = 1375
$ws.Range("I22").Value = 1393.75
$ws.Range("J22").Value = 1300
$ws.Range("K22").Value = 1393.75
$ws.Range("L22").Value = 1300
$ws.Range("M22").Value = -1098.75
$ws.Range("N22").Value = -1890
$ws.Range("H27").Value = 1375
$ws.Range("I27").Value = 1393.75
$ws.Range("J27").Value = 1300
$ws.Range("K27").Value = 1393.75
$ws.Range("L27").Value = 1300
$ws.Range("M27").Value = -1286.75
$ws.Range("N27").Value = -1514
$ws.Range("H61").Value = 1950
$ws.Range("I61").Value = 1950
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1950
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1748
$ws.Range("H113").Value = 1950
$ws.Range("I113").Value = 1950
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1950
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 220
$ws.Range("H132").Value = 8593.23
$ws.Range("I132").Value = 8453
$ws.Range("J132").Value = 8713.429
$ws.Range("K132").Value = 25359
$ws.Range("L132").Value = 26140.287
$ws.Range("M132").Value = -22829
$ws.Range("N132").Value = -31200.287

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 1499.5
$ws.Range("K132").Value = 4500
$ws.Range("L132").Value = 4498.5
$ws.Range("M132").Value = -1970
$ws.Range("N132").Value = -9558.5
